# Generate Report for Handback
#
# The localization-status report records, per target language, the
# handback state of each source file. The file
# f2a91ce9-c14c-402b-8b0d-615309079abd.md has now been handed back in
# sync with en-US, so its Status flips from "Ready for handoff" to
# "Handed back: in sync with en-US" everywhere it is reported, and the
# per-language "Latest Handback DateTime" is stamped with the handback
# time for each locale.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: one summary row per language for this file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# --- zh-cn detail sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = $newStatus
$wsZhCn.Range("G3").Value = "2016-03-02 10:39:06"

# --- de-de detail sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = $newStatus
$wsDeDe.Range("G3").Value = "2016-03-02 10:39:25"
